# Switch the single section's page setup from portrait to landscape.
# Word automatically swaps the page width/height when the orientation
# changes, turning the 11952 x 16848 (w x h, twips) portrait page into
# a 16848 x 11952 landscape page.
$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    $sec.PageSetup.Orientation = 1  # wdOrientLandscape
}
